# Bharadwaj_LabExam03Grading.xlsx - "adloori to davuluri completed"
# Fill in the "Total Points" (grading) column E for the Customer Class and
# Product Class sections, which were previously left blank while grading
# was in progress. The section/grand totals recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Customer Class section (rows 3-6)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Product Class section (rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Reset the scroll position back to the top-left (the saved view had
# scrolled down to row 20) and move the active selection to E15.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E15").Select()
